# Apply the "automatic update" edit to the Artfynd sheet.
# Rows 76-79 have most of their field values rotated (row76<-old77,
# row77<-old78, row78<-old79, row79<-old76), while column B
# (Taxonsorteringsordning) gets its own independent set of new values.
# Rows 75 and 81 only get their B value updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75: only the taxon sort order (B) changes ---
$ws.Range("B75").Value = 96735

# --- Row 76 (now holds what used to be row 77's record) ---
$ws.Range("A76").Value = 111950184
$ws.Range("B76").Value = 56575
$ws.Range("D76").Value = "NT"
$ws.Range("E76").Value = 103021
$ws.Range("F76").Value = "Talltita"
$ws.Range("G76").Value = "Poecile montanus"
$ws.Range("H76").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I76").ClearContents()
$ws.Range("J76").ClearContents()
$ws.Range("Q76").Value = 580447
$ws.Range("R76").Value = 7053302
$ws.Range("S76").Value = 10
$ws.Range("Z76").Value = "19:37"
$ws.Range("AB76").Value = "19:37"
$ws.Range("AW76").Value = "Kim Hultgren"
$ws.Range("AX76").Value = "Kim Hultgren"

# --- Row 77 (now holds what used to be row 78's record) ---
$ws.Range("A77").Value = 111949317
$ws.Range("B77").Value = 96652
$ws.Range("D77").Value = "LC"
$ws.Range("E77").Value = 219790
$ws.Range("F77").Value = "Fläcknycklar"
$ws.Range("G77").Value = "Dactylorhiza maculata"
$ws.Range("H77").Value = "(L.) Soó"
$ws.Range("Q77").Value = 580500
$ws.Range("R77").Value = 7053329
$ws.Range("S77").Value = 2
$ws.Range("Z77").Value = "18:54"
$ws.Range("AB77").Value = "18:54"

# --- Row 78 (now holds what used to be row 79's record) ---
$ws.Range("A78").Value = 111949678
$ws.Range("B78").Value = 96735
$ws.Range("D78").Value = "VU"
$ws.Range("E78").Value = 220787
$ws.Range("F78").Value = "Knärot"
$ws.Range("G78").Value = "Goodyera repens"
$ws.Range("H78").Value = "(L.) R. Br."
# Column I ("Antal") holds numeric-looking values but must stay TEXT,
# matching the workbook's original inline-string typing. Force a text
# number format before assigning so Excel doesn't coerce it to a number.
$ws.Range("I78").NumberFormat = "@"
$ws.Range("I78").Value = "7"
$ws.Range("Q78").Value = 580467
$ws.Range("R78").Value = 7053330
$ws.Range("Z78").Value = "19:11"
$ws.Range("AB78").Value = "19:11"

# --- Row 79 (now holds what used to be row 76's record) ---
$ws.Range("A79").Value = 111949575
$ws.Range("B79").Value = 96735
$ws.Range("I79").NumberFormat = "@"
$ws.Range("I79").Value = "15"
$ws.Range("J79").Value = "plantor/tuvor"
$ws.Range("Q79").Value = 580471
$ws.Range("R79").Value = 7053333
$ws.Range("S79").Value = 1
$ws.Range("Z79").Value = "19:05"
$ws.Range("AB79").Value = "19:05"
$ws.Range("AW79").Value = "Kamilla Andersson"
$ws.Range("AX79").Value = "Kamilla Andersson"

# --- Row 81: only the taxon sort order (B) changes ---
$ws.Range("B81").Value = 90808
